$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(32, 1).Value = "Gretchen Flores Calica"
$ws.Cells.Item(32, 2).Value = "女"
$ws.Cells.Item(32, 3).Value = "Teach International (国际认证)"
$ws.Cells.Item(32, 4).Value = "TESOL & TEYL (140小时双证)"
$ws.Cells.Item(32, 5).Value = "2019年获证 (资深)"
$ws.Cells.Item(32, 6).Value = "英语母语/国际外教"
$ws.Cells.Item(32, 7).Value = "青少年英语(TEYL) | 趣味互动 | 语法"
$ws.Cells.Item(32, 8).Value = "1.持有140小时高阶证书(含TEYL青少年专项)。`n2.专攻青少年心理发展阶段，懂孩子。`n3.擅长利用游戏和趣味活动活跃课堂。"
$ws.Cells.Item(32, 9).Value = "1.接受过系统的“教师工具箱”培训，资源丰富。`n2.擅长课程规划与时间管理。`n3.专注于EFL/ESL(非母语英语)环境教学。"

$ws.Cells.Item(33, 1).Value = "Jackylou Mariano"
$ws.Cells.Item(33, 2).Value = "女"
$ws.Cells.Item(33, 3).Value = "World TESOL Academy (英国认证)"
$ws.Cells.Item(33, 4).Value = "TESOL/TEFL (120小时)"
$ws.Cells.Item(33, 5).Value = "2021年获证"
$ws.Cells.Item(33, 6).Value = "英语母语/国际外教"
$ws.Cells.Item(33, 7).Value = "在线英语教学 | 海外EFL | 口语"
$ws.Cells.Item(33, 8).Value = "1.英国UKRLP注册机构认证，资质过硬。`n2.经过专门的“在线课堂环境”教学培训。`n3.教学风格标准化，适应性强。"
$ws.Cells.Item(33, 9).Value = "1.拥有ACCREDITAT国际认证。`n2.专注于非母语学生的海外及在线教学。`n3.具备处理跨文化交流的经验。"

$ws.Cells.Item(34, 1).Value = "Nikki Elaine Pertubal"
$ws.Cells.Item(34, 2).Value = "女"
$ws.Cells.Item(34, 3).Value = "TEFL Pro Institute (国际专业发展学院)"
$ws.Cells.Item(34, 4).Value = "TEFL (120小时)"
$ws.Cells.Item(34, 5).Value = "2021年获证"
$ws.Cells.Item(34, 6).Value = "英语母语/国际外教"
$ws.Cells.Item(34, 7).Value = "全年龄段英语 | 听说读写 | 词汇"
$ws.Cells.Item(34, 8).Value = "1.核心单元全优通过(Distinction)。`n2.擅长课堂管理，能搞定坐不住的孩子。`n3.覆盖从青少年到成人的全体系教学能力。"
$ws.Cells.Item(34, 9).Value = "1.系统完成了120小时核心单元训练。`n2.擅长根据学生反馈调整教学策略。`n3.拥有扎实的词汇与发音教学实战技巧。"

$ws.Cells.Item(35, 1).Value = "Mark Diaz"
$ws.Cells.Item(35, 2).Value = "男"
$ws.Cells.Item(35, 3).Value = "TEFL Pro Institute (国际专业发展学院)"
$ws.Cells.Item(35, 4).Value = "TEFL (120小时)"
$ws.Cells.Item(35, 5).Value = "2021年获证"
$ws.Cells.Item(35, 6).Value = "英语母语/国际外教"
$ws.Cells.Item(35, 7).Value = "成人/青少年英语 | 发音纠正 | 写作"
$ws.Cells.Item(35, 8).Value = "1.擅长教授听觉语言技能(说与听)。`n2.对读写技能(读与写)有系统教学法。`n3.善于处理当代英语语言学习中的难点。"
$ws.Cells.Item(35, 9).Value = "1.完成17个核心教学模块的系统训练。`n2.擅长针对青少年和成人不同心理特点备课。`n3.注重培养学生的实际语言运用能力。"

$ws.Cells.Item(36, 1).Value = "Dennis Arvin Junatas"
$ws.Cells.Item(36, 2).Value = "男"
$ws.Cells.Item(36, 3).Value = "TEFL Pro Institute (国际专业发展学院)"
$ws.Cells.Item(36, 4).Value = "TEFL (120小时)"
$ws.Cells.Item(36, 5).Value = "2021年获证"
$ws.Cells.Item(36, 6).Value = "英语母语/国际外教"
$ws.Cells.Item(36, 7).Value = "基础英语巩固 | 课堂互动 | 语法"
$ws.Cells.Item(36, 8).Value = "1.擅长学习者反馈，能敏锐发现学生问题。`n2.精通语法与发音教学。`n3.具备扎实的对外英语教学理论基础。"
$ws.Cells.Item(36, 9).Value = "1.熟练掌握课堂管理技巧，课堂氛围好。`n2.对英语语法教学有独到的拆解方法。`n3.擅长引导学生开口说英语。"

$ws.Cells.Item(37, 1).Value = "Myla Lalaine B. Uchu-e"
$ws.Cells.Item(37, 2).Value = "女"
$ws.Cells.Item(37, 3).Value = "TEFL Pro Institute (国际专业发展学院)"
$ws.Cells.Item(37, 4).Value = "TEFL (120小时)"
$ws.Cells.Item(37, 5).Value = "持证外教"
$ws.Cells.Item(37, 6).Value = "英语母语/国际外教"
$ws.Cells.Item(37, 7).Value = "综合英语 | 学习策略 | 国际教学"
$ws.Cells.Item(37, 8).Value = "1.以优异成绩(High Distinction)通过考核。`n2.精通不同教学方法在多样化需求中的应用。`n3.擅长用易于理解的方式讲解复杂语法。"
$ws.Cells.Item(37, 9).Value = "1.具备国际教学知识(International Teaching Knowledge)。`n2.课堂管理能力强，教学逻辑清晰。`n3.能根据学生需求定制学习策略。"
